# LM5060 calculation update: refresh resistor/parameter values and
# restore the pin-map selection on the Datasheet sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datasheet")

# --- Solve & Design section ---
# VINMAX(V)
$ws.Range("C9").Value = 8
# R1(kOhm)
$ws.Range("C10").Value = 26
# R2(kOhm)
$ws.Range("C11").Value = 66.5

# --- Pick & Evaluate section ---
# R2(kOhm)
$ws.Range("C18").Value = 11
# R3(kOhm)
$ws.Range("C19").Value = 6.49

# Restore the active selection / pin map to H19
$ws.Activate()
$ws.Range("H19").Select()
